# Auto-generated edit script applying cached-value updates to Sheets/Spriggan_Profits.xlsx
# (workbook sheet tab names: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 517.2222
$ws.Range("I6").Value = 269.375
$ws.Range("K6").Value = 808.125
$ws.Range("M6").Value = -696.125
$ws.Range("H17").Value = 486988.4
$ws.Range("J17").Value = 525787.75
$ws.Range("L17").Value = 1577363.25
$ws.Range("N17").Value = -1577699.25
$ws.Range("H42").Value = 2920.2
$ws.Range("I42").Value = 2900.5
$ws.Range("J42").Value = 2999
$ws.Range("K42").Value = 8701.5
$ws.Range("L42").Value = 8997
$ws.Range("M42").Value = -8471.5
$ws.Range("N42").Value = -9457
$ws.Range("H80").Value = 500.75
$ws.Range("I80").Value = 549
$ws.Range("J80").Value = 484.66666
$ws.Range("K80").Value = 1647
$ws.Range("L80").Value = 1453.99998
$ws.Range("M80").Value = -649
$ws.Range("N80").Value = -3449.99998
$ws.Range("H83").Value = 500.75
$ws.Range("I83").Value = 549
$ws.Range("J83").Value = 484.66666
$ws.Range("K83").Value = 4941
$ws.Range("L83").Value = 4361.99994
$ws.Range("M83").Value = 51
$ws.Range("N83").Value = -14345.99994
$ws.Range("H86").Value = 9594.916999999999
$ws.Range("I86").Value = 10103.546
$ws.Range("K86").Value = 10103.546
$ws.Range("M86").Value = -8980.546
$ws.Range("H89").Value = 9594.916999999999
$ws.Range("I89").Value = 10103.546
$ws.Range("K89").Value = 50517.73
$ws.Range("M89").Value = -44901.73
$ws.Range("H132").Value = 1477.0834
$ws.Range("I132").Value = 1474.2188
$ws.Range("K132").Value = 4422.6564
$ws.Range("M132").Value = -1892.6564
$ws.Range("H138").Value = 5847.3735
$ws.Range("I138").Value = 4588.5625
$ws.Range("J138").Value = 6147.985
$ws.Range("K138").Value = 13765.6875
$ws.Range("L138").Value = 18443.955
$ws.Range("M138").Value = -8625.6875
$ws.Range("N138").Value = -28723.955

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1646.5151
$ws.Range("I2").Value = 1480.8
$ws.Range("J2").Value = 2164.375
$ws.Range("K2").Value = 1480.8
$ws.Range("L2").Value = 2164.375
$ws.Range("M2").Value = -1367.8
$ws.Range("N2").Value = -2390.375
$ws.Range("H34").Value = 39797.2
$ws.Range("I34").Value = 38990
$ws.Range("K34").Value = 38990
$ws.Range("M34").Value = -38719
$ws.Range("H61").Value = 200002190
$ws.Range("I61").Value = 250001500
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 250001500
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -250001288
$ws.Range("N61").Value = -5424
$ws.Range("H95").Value = 40867.832
$ws.Range("J95").Value = 40867.832
$ws.Range("L95").Value = 40867.832
$ws.Range("N95").Value = -46359.832
$ws.Range("H116").Value = 1646.5151
$ws.Range("I116").Value = 1480.8
$ws.Range("J116").Value = 2164.375
$ws.Range("K116").Value = 1480.8
$ws.Range("L116").Value = 2164.375
$ws.Range("M116").Value = 813.2
$ws.Range("N116").Value = -6752.375
$ws.Range("H132").Value = 4527530
$ws.Range("I132").Value = 2442835.5
$ws.Range("K132").Value = 7328506.5
$ws.Range("M132").Value = -7325976.5
$ws.Range("H136").Value = 200002190
$ws.Range("I136").Value = 250001500
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 750004500
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -750001950
$ws.Range("N136").Value = -20100

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1646.5151
$ws.Range("I3").Value = 1480.8
$ws.Range("J3").Value = 2164.375
$ws.Range("K3").Value = 1480.8
$ws.Range("L3").Value = 2164.375
$ws.Range("M3").Value = -1366.8
$ws.Range("N3").Value = -2392.375
$ws.Range("H103").Value = 15936.75
$ws.Range("J103").Value = 15936.75
$ws.Range("L103").Value = 15936.75
$ws.Range("N103").Value = -18280.75
$ws.Range("H105").Value = 2963.9333
$ws.Range("I105").Value = 1962.375
$ws.Range("K105").Value = 1962.375
$ws.Range("M105").Value = -215.375
$ws.Range("H107").Value = 201494
$ws.Range("I107").Value = 1867.75
$ws.Range("J107").Value = 999999
$ws.Range("K107").Value = 1867.75
$ws.Range("L107").Value = 999999
$ws.Range("M107").Value = 52.25
$ws.Range("N107").Value = -1003839
$ws.Range("H109").Value = 66666
$ws.Range("J109").Value = 66666
$ws.Range("L109").Value = 66666
$ws.Range("N109").Value = -69440
$ws.Range("H130").Value = 94495
$ws.Range("J130").Value = 94495
$ws.Range("L130").Value = 94495
$ws.Range("N130").Value = -104535

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 95
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 95
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 95
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -321
$ws.Range("H10").Value = 422.5
$ws.Range("I10").Value = 399
$ws.Range("J10").Value = 446
$ws.Range("K10").Value = 399
$ws.Range("L10").Value = 446
$ws.Range("M10").Value = -260
$ws.Range("N10").Value = -724
$ws.Range("H132").Value = 90911650
$ws.Range("I132").Value = 100002310
$ws.Range("K132").Value = 300006930
$ws.Range("M132").Value = -300004400

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9900.5
$ws.Range("I3").Value = 9900.5
$ws.Range("K3").Value = 29701.5
$ws.Range("M3").Value = -29589.5
$ws.Range("H12").Value = 412
$ws.Range("I12").Value = 818.25
$ws.Range("K12").Value = 2454.75
$ws.Range("M12").Value = -2281.75
$ws.Range("H22").Value = 242.85715
$ws.Range("I22").Value = 241.66667
$ws.Range("J22").Value = 250
$ws.Range("K22").Value = 725.00001
$ws.Range("L22").Value = 750
$ws.Range("M22").Value = -556.00001
$ws.Range("N22").Value = -1088
$ws.Range("H27").Value = 242.85715
$ws.Range("I27").Value = 241.66667
$ws.Range("J27").Value = 250
$ws.Range("K27").Value = 725.00001
$ws.Range("L27").Value = 750
$ws.Range("M27").Value = -623.00001
$ws.Range("N27").Value = -954
$ws.Range("H117").Value = 3015
$ws.Range("J117").Value = 3000
$ws.Range("L117").Value = 9000
$ws.Range("N117").Value = -15884
$ws.Range("H132").Value = 1998
$ws.Range("J132").Value = 1998
$ws.Range("L132").Value = 17982
$ws.Range("N132").Value = -23042

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H80").Value = 3129.4
$ws.Range("I80").Value = 3237
$ws.Range("K80").Value = 3237
$ws.Range("M80").Value = -2239
$ws.Range("H83").Value = 3129.4
$ws.Range("I83").Value = 3237
$ws.Range("K83").Value = 16185
$ws.Range("M83").Value = -11193
$ws.Range("H102").Value = 4922.1714
$ws.Range("I102").Value = 2186.7812
$ws.Range("K102").Value = 2186.7812
$ws.Range("M102").Value = -564.7811999999999
$ws.Range("H113").Value = 33021.625
$ws.Range("I113").Value = 39605.188
$ws.Range("J113").Value = 6687.375
$ws.Range("K113").Value = 39605.188
$ws.Range("L113").Value = 6687.375
$ws.Range("M113").Value = -37435.188
$ws.Range("N113").Value = -11027.375
$ws.Range("H132").Value = 2286884.5
$ws.Range("I132").Value = 2328308.2
$ws.Range("K132").Value = 6984924.600000001
$ws.Range("M132").Value = -6982394.600000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4998.5
$ws.Range("I122").Value = 4998.6665
$ws.Range("K122").Value = 14995.9995
$ws.Range("M122").Value = -12545.9995
$ws.Range("H132").Value = 9266102
$ws.Range("I132").Value = 11370843
$ws.Range("J132").Value = 5239.8
$ws.Range("K132").Value = 34112529
$ws.Range("L132").Value = 15719.4
$ws.Range("M132").Value = -34109999
$ws.Range("N132").Value = -20779.4
$ws.Range("H136").Value = 2858.7896
$ws.Range("I136").Value = 2765.818
$ws.Range("K136").Value = 8297.454000000002
$ws.Range("M136").Value = -5747.454000000002

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 19500
$ws.Range("I2").Value = 19500
$ws.Range("K2").Value = 19500
$ws.Range("M2").Value = -19388
$ws.Range("H24").Value = 20006
$ws.Range("I24").Value = 9
$ws.Range("K24").Value = 9
$ws.Range("M24").Value = 221
$ws.Range("H62").Value = 5965.5557
$ws.Range("J62").Value = 6311.25
$ws.Range("L62").Value = 6311.25
$ws.Range("N62").Value = -7559.25
$ws.Range("H65").Value = 5965.5557
$ws.Range("J65").Value = 6311.25
$ws.Range("L65").Value = 31556.25
$ws.Range("N65").Value = -37796.25
$ws.Range("H107").Value = 1106.4517
$ws.Range("I107").Value = 655.25
$ws.Range("J107").Value = 1926.8182
$ws.Range("K107").Value = 1965.75
$ws.Range("L107").Value = 5780.4546
$ws.Range("M107").Value = -45.75
$ws.Range("N107").Value = -9620.454600000001
$ws.Range("H122").Value = 3412.4285
$ws.Range("I122").Value = 3581.5
$ws.Range("K122").Value = 10744.5
$ws.Range("M122").Value = -8294.5
$ws.Range("H136").Value = 19232922
$ws.Range("I136").Value = 20002098
$ws.Range("K136").Value = 60006294
$ws.Range("M136").Value = -60003744
